# Update 2p0. Convention change to support multi-axle vehicles
#
# Adds two new vehicle-body CAD sheets ("Truck_Amandla" and
# "Trailer_Kumanzi"), positioned either side of the existing
# "Trailer_Thwala" sheet, using "Bus_Makhulu" (an existing 8-row sheet
# with the same layout/styles) as the structural template.
#
# NOTE: worksheet variables captured from $wb.Worksheets.Item(...) track
# the sheet *position*, not sheet identity -- once another sheet is
# inserted/moved and shifts positions around, an old variable can now
# silently resolve to a different sheet. To stay safe we re-fetch each
# sheet by its (now final) Name right before editing its cells, and we
# finish all moves before doing any value edits.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Bus_Makhulu")

# --- Truck_Amandla: create as a copy of Bus_Makhulu, place immediately
#     before Trailer_Thwala ---
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "Truck_Amandla"
$wb.Worksheets.Item("Truck_Amandla").Move($wb.Worksheets.Item("Trailer_Thwala"))

# --- Trailer_Kumanzi: create as a copy of Bus_Makhulu, place immediately
#     after Trailer_Thwala ---
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "Trailer_Kumanzi"
$wb.Worksheets.Item("Trailer_Kumanzi").Move($null, $wb.Worksheets.Item("Trailer_Thwala"))

# --- Fill in Truck_Amandla's data (re-fetched by name) ---
$truck = $wb.Worksheets.Item("Truck_Amandla")

$truck.Range("H3").Value = "CAD_Truck_Amandla"
$truck.Range("H4").Value = "CAD_Truck_Amandla"

$truck.Range("F5").Value = 0
$truck.Range("G5").Value = 0
$truck.Range("H5").Value = 0

$truck.Range("F6").Value = 0
$truck.Range("G6").Value = 0
$truck.Range("H6").Value = 0

$truck.Range("F7").Value = 0.6
$truck.Range("G7").Value = 0.8
$truck.Range("H7").Value = 1

$truck.Range("H8").Value = 1

# --- Fill in Trailer_Kumanzi's data (re-fetched by name) ---
$kumanzi = $wb.Worksheets.Item("Trailer_Kumanzi")

$kumanzi.Range("H3").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "CAD_Trailer_Kumanzi"

$kumanzi.Range("F5").Value = 0
$kumanzi.Range("G5").Value = 0
$kumanzi.Range("H5").Value = 0

$kumanzi.Range("F6").Value = 0
$kumanzi.Range("G6").Value = 0
$kumanzi.Range("H6").Value = 0

$kumanzi.Range("F7").Value = 1
$kumanzi.Range("G7").Value = 0.75
$kumanzi.Range("H7").Value = 0.055

$kumanzi.Range("H8").Value = 0.5

# Trailer_Kumanzi ends up the active/selected sheet.
$wb.Worksheets.Item("Trailer_Kumanzi").Activate()
